# Edit script: rewrite "Enigmatic Symmetries" essay into
# "The Symphony of the Human Body" essay, update author/email, and
# append a trailing empty paragraph.
#
# We locate each run's old text with Find (no replace, so the engine's
# smart-quote/auto-format substitution never kicks in) and then assign
# the new literal text straight onto the found Range.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                                $true, 1, $false)
    if (-not $found) {
        Write-Host "NOT FOUND:" $old
        return
    }
    $rng.Text = $new
}

# --- Title ---------------------------------------------------------
Replace-Text "Enigmatic Symmetries in Nature's Symphony" `
             "The Symphony of the Human Body: Exploring the Wonders of Life"

# --- Author line -----------------------------------------------------
Replace-Text " Stella Maris" " Clara Richards"

# --- Email line (merges "stella" / "." / "maris@cosmicmelodies" / "." / "edu") ---
Replace-Text "stella" "clararichards@edumail"
Replace-Text "maris@cosmicmelodies.edu" "org"

# --- Body paragraph --------------------------------------------------
Replace-Text "In the grand symphony of existence, nature orchestrates intricate patterns governed by enigmatic symmetries" `
             "In the vast expanse of the cosmos, there lies a microcosm of wonders--the human body"

Replace-Text " From the subatomic realm to the expansive cosmos, these symmetries reveal an elegance and harmony beyond human comprehension" `
             " A captivating symphony of interconnected systems, each working harmoniously to sustain life, this intricate machinery holds boundless mysteries waiting to be unraveled"

Replace-Text " In this cosmic ballet, fundamental particles waltz in harmonious synchronization, their interactions governed by the ethereal laws of quantum mechanics" `
             " From the intricate choreography of cells to the symphony of organs, the human body is a testament to the marvels of nature"

# Runs " The dance of molecules...world" + "." + " Celestial bodies...forces" collapse into one run.
Replace-Text " The dance of molecules, guided by chemical affinities, sculpts the diverse tapestry of substances that grace our world. Celestial bodies, from planets to galaxies, trace out graceful arcs in a cosmic ballet choreographed by gravitational forces" `
             " It's a dynamic masterpiece that deserves our profound admiration and exploration"

Replace-Text "As we ascend the ladder of complexity, biological systems unveil a symphony of symmetries, from the intricate arrangements of DNA to the mesmerizing patterns of animal behavior" `
             "Enter the world of biology, a field that unlocks the secrets of life's mechanisms"

Replace-Text " Evolution, the maestro of life, has crafted organisms that possess both symmetry and asymmetry, each serving a vital purpose in nature's grand design" `
             " It unravels the enigmatic blueprint of DNA, the blueprint of life, guiding the development and functioning of every living organism"

Replace-Text " Even in the seemingly chaotic realm of human societies, patterns and symmetries emerge, shaped by cultural norms, economic forces, and political structures" `
             " Biology illuminates the intricacies of cellular processes, revealing the hidden language of molecules and their profound impact on our being"

Replace-Text " These symmetries, both tangible and abstract, provide a glimpse into the underlying order that permeates all aspects of existence" `
             " It explores the marvelous tapestry of ecosystems, highlighting the interconnectedness of all living creatures, and unveils the evolutionary saga that has shaped the diversity of life on Earth"

Replace-Text "Exploring these symmetries not only unravels the mysteries of nature but also offers practical benefits" `
             "Within the vast canvas of biology, the human body stands as a captivating subject of study"

Replace-Text " The insights gained from studying symmetries have led to groundbreaking advancements in physics, chemistry, biology, and engineering" `
             " Its intricate symphony of organs and tissues, working in harmony, enables us to experience the world around us"

# Runs " Symmetry considerations...lives" + "." + " Delving into...universe" collapse into one run.
Replace-Text " Symmetry considerations have guided the design of new materials, drugs, and technologies that have revolutionized our lives. Delving into the enigmatic world of symmetries is akin to embarking on a quest for hidden treasures, revealing the beauty and interconnectedness of the universe" `
             " Biology allows us to delve into the mechanisms of digestion, respiration, and circulation, understanding how our bodies transform nutrients into energy, breathe life-giving oxygen, and circulate blood throughout our intricate network of vessels"

# --- Summary paragraph -------------------------------------------------
Replace-Text "The exploration of symmetries in nature unveils a breathtaking tapestry of patterns and harmonies that span the entire spectrum of existence" `
             "The human body is an enigma, a mesmerizing symphony of intricate systems that orchestrates the miracle of life"

Replace-Text " From the subatomic realm to the vast reaches of the cosmos, symmetries reveal an underlying elegance and order that govern the universe" `
             " Biology, a field dedicated to unraveling the tapestry of life, shines a light on the inner workings of our bodies, unveiling the secrets of our existence"

# Runs " The study...the" + "fundamental laws...technology" + "." + " Unraveling...humanity" collapse into one run.
Replace-Text " The study of these symmetries has yielded profound insights into the fundamental laws of nature and has led to transformative advancements in various fields of science and technology. Unraveling the enigmas of symmetries not only enriches our understanding of the universe but also empowers us to harness its mysteries for the betterment of humanity" `
             " From cellular processes to organ functions, the study of biology illuminates the wonder of life and inspires us to delve deeper into the complexities of the human experience"

# --- Append a trailing empty paragraph at the end of the document ------
[void]$d.Paragraphs.Add()
